$wb = $excel.ActiveWorkbook

# Rename existing sheet "Hoja1" -> "Hitos"
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Name = "Hitos"

# Add a new worksheet for "Tareas divididas" after "Hitos"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tareas divididas"

# Fill in new sheet data.
# Shared-string pool order matters (new strings are appended in first-use
# order): Leandro, Ezequiel, Franco, then the 3 task descriptions in
# B3, B1, B2 order (ABM Productos.., ABM Usuarios.., Funcion de mesero).
$ws2.Range("A1").Value = "Leandro"
$ws2.Range("A2").Value = "Ezequiel"
$ws2.Range("A3").Value = "Franco"
$ws2.Range("B3").Value = "ABM de Productos, Categorias, Sub Categorias"
$ws2.Range("B1").Value = "ABM de Usuarios, Listado de mesas"
$ws2.Range("B2").Value = "Funcion de mesero"

# Set column widths on new sheet (target stored widths: 18 and ~42.285).
# The host's ColumnWidth setter re-quantizes to whole-pixel units (1/6 char
# here), so feed it the pre-image that lands on the desired stored width.
$ws2.Columns.Item(1).ColumnWidth = 17.166666666666668
$ws2.Columns.Item(2).ColumnWidth = 41.451822916666664

# Select cell on new sheet (as in diff: activeCell B4)
$ws2.Range("B4").Select()

# Select cell on Hitos sheet (as in diff: activeCell C5)
$ws1.Range("C5").Select()

# Make "Tareas divididas" the active sheet (tabSelected / activeTab=1)
$ws2.Activate()
